$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 14671.286
$ws.Range("I12").Value = 20289.4
$ws.Range("J12").Value = 626
$ws.Range("K12").Value = 20289.4
$ws.Range("L12").Value = 626
$ws.Range("M12").Value = -20119.4
$ws.Range("N12").Value = -966

$ws.Range("H38").Value = 1379.909
$ws.Range("I38").Value = 1379.909
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 4139.727000000001
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = ""
$ws.Range("N38").Value = -3767.727000000001

$ws.Range("H42").Value = 4405.6
$ws.Range("I42").Value = 3673
$ws.Range("K42").Value = 11019
$ws.Range("M42").Value = -10789

$ws.Range("H48").Value = 9115.4
$ws.Range("J48").Value = 9398.5
$ws.Range("L48").Value = 28195.5
$ws.Range("N48").Value = -28779.5

$ws.Range("H51").Value = 10437.125
$ws.Range("I51").Value = 3249
$ws.Range("K51").Value = 3249
$ws.Range("M51").Value = -2765

$ws.Range("H56").Value = 9115.4
$ws.Range("J56").Value = 9398.5
$ws.Range("L56").Value = 28195.5
$ws.Range("N56").Value = -29263.5

$ws.Range("H74").Value = 14527163
$ws.Range("I74").Value = 14527163
$ws.Range("K74").Value = 14527163
$ws.Range("M74").Value = -14526227

$ws.Range("H77").Value = 14527163
$ws.Range("I77").Value = 14527163
$ws.Range("K77").Value = 72635815
$ws.Range("M77").Value = -72631135

$ws.Range("H86").Value = 2380.3333
$ws.Range("I86").Value = 2570.5
$ws.Range("K86").Value = 2570.5
$ws.Range("M86").Value = -1447.5

$ws.Range("H89").Value = 2380.3333
$ws.Range("I89").Value = 2570.5
$ws.Range("K89").Value = 12852.5
$ws.Range("M89").Value = -7236.5

$ws.Range("H103").Value = 2901.4167
$ws.Range("I103").Value = 799.3333
$ws.Range("J103").Value = 3602.111
$ws.Range("K103").Value = 2397.9999
$ws.Range("L103").Value = 10806.333
$ws.Range("M103").Value = -1811.9999
$ws.Range("N103").Value = -11978.333

$ws.Range("H112").Value = 67033.61
$ws.Range("J112").Value = 55471.633
$ws.Range("L112").Value = 166414.899
$ws.Range("N112").Value = -168630.899

$ws.Range("H138").Value = 3880.5
$ws.Range("J138").Value = 5136.3335
$ws.Range("L138").Value = 15409.0005
$ws.Range("N138").Value = -25689.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 17379
$ws.Range("J46").Value = 24999.5
$ws.Range("L46").Value = 24999.5
$ws.Range("N46").Value = -25637.5

$ws.Range("H63").Value = 4593
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 4593
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = ""
$ws.Range("M63").Value = 4593
$ws.Range("N63").Value = -5965

$ws.Range("H66").Value = 4593
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 4593
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = ""
$ws.Range("M66").Value = 22965
$ws.Range("N66").Value = -29829

$ws.Range("H74").Value = 55561050
$ws.Range("I74").Value = 62505744
$ws.Range("K74").Value = 62505744
$ws.Range("M74").Value = -62504870

$ws.Range("H77").Value = 55561050
$ws.Range("I77").Value = 62505744
$ws.Range("K77").Value = 312528720
$ws.Range("M77").Value = -312524352

$ws.Range("H88").Value = 1932.0454
$ws.Range("I88").Value = 1687.7142
$ws.Range("J88").Value = 2046.0667
$ws.Range("K88").Value = 1687.7142
$ws.Range("L88").Value = 2046.0667
$ws.Range("M88").Value = -1281.7142
$ws.Range("N88").Value = -2858.0667

$ws.Range("H91").Value = 1932.0454
$ws.Range("I91").Value = 1687.7142
$ws.Range("J91").Value = 2046.0667
$ws.Range("K91").Value = 1687.7142
$ws.Range("L91").Value = 2046.0667
$ws.Range("M91").Value = -283.7141999999999
$ws.Range("N91").Value = -4854.0667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5004.16
$ws.Range("I86").Value = 4560.9
$ws.Range("J86").Value = 6777.2
$ws.Range("K86").Value = 4560.9
$ws.Range("L86").Value = 6777.2
$ws.Range("M86").Value = -3437.9
$ws.Range("N86").Value = -9023.200000000001

$ws.Range("H89").Value = 5004.16
$ws.Range("I89").Value = 4560.9
$ws.Range("J89").Value = 6777.2
$ws.Range("K89").Value = 22804.5
$ws.Range("L89").Value = 33886
$ws.Range("M89").Value = -17188.5
$ws.Range("N89").Value = -45118

$ws.Range("H118").Value = 125706.664
$ws.Range("J118").Value = 125706.664
$ws.Range("L118").Value = 125706.664
$ws.Range("N118").Value = -129020.664

$ws.Range("H134").Value = 12823289
$ws.Range("I134").Value = 14707894
$ws.Range("K134").Value = 44123682
$ws.Range("M134").Value = -44121147

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6531.778
$ws.Range("I62").Value = 3830
$ws.Range("J62").Value = 7882.6665
$ws.Range("K62").Value = 3830
$ws.Range("L62").Value = 7882.6665
$ws.Range("M62").Value = -3206
$ws.Range("N62").Value = -9130.666499999999

$ws.Range("H65").Value = 6531.778
$ws.Range("I65").Value = 3830
$ws.Range("J65").Value = 7882.6665
$ws.Range("K65").Value = 19150
$ws.Range("L65").Value = 39413.3325
$ws.Range("M65").Value = -16030
$ws.Range("N65").Value = -45653.3325

$ws.Range("H107").Value = 496316.9
$ws.Range("I107").Value = 1553342
$ws.Range("K107").Value = 1553342
$ws.Range("M107").Value = -1551422

$ws.Range("H134").Value = 10418844
$ws.Range("I134").Value = 12501942
$ws.Range("K134").Value = 37505826
$ws.Range("M134").Value = -37503291

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1390
$ws.Range("I3").Value = 275
$ws.Range("J3").Value = 3620
$ws.Range("K3").Value = 825
$ws.Range("L3").Value = 10860
$ws.Range("M3").Value = -713
$ws.Range("N3").Value = -11084

$ws.Range("H11").Value = 121040
$ws.Range("I11").Value = 125666.664
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 376999.992
$ws.Range("L11").Value = 30000
$ws.Range("M11").Value = -376859.992
$ws.Range("N11").Value = -30280

$ws.Range("H113").Value = 100479.7
$ws.Range("I113").Value = 250411
$ws.Range("K113").Value = 751233
$ws.Range("M113").Value = -749063

$ws.Range("H133").Value = 11773.357
$ws.Range("I133").Value = 8482.799999999999
$ws.Range("K133").Value = 25448.4
$ws.Range("M133").Value = -20388.4

$ws.Range("H134").Value = 1511
$ws.Range("I134").Value = 1511
$ws.Range("K134").Value = 4533
$ws.Range("M134").Value = 537

$ws.Range("H137").Value = 6251403.5
$ws.Range("I137").Value = 7693600.5
$ws.Range("J137").Value = 1883
$ws.Range("K137").Value = 23080801.5
$ws.Range("L137").Value = 5649
$ws.Range("M137").Value = -23075701.5
$ws.Range("N137").Value = -15849

$ws.Range("H138").Value = 2372.6667
$ws.Range("I138").Value = 2531.889
$ws.Range("K138").Value = 7595.667
$ws.Range("M138").Value = -2455.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 29666.334
$ws.Range("I10").Value = 29500
$ws.Range("K10").Value = 29500
$ws.Range("M10").Value = -29331

$ws.Range("H18").Value = 39666.668
$ws.Range("I18").Value = 42000
$ws.Range("K18").Value = 42000
$ws.Range("M18").Value = -41707

$ws.Range("H80").Value = 3360
$ws.Range("I80").Value = 2950
$ws.Range("K80").Value = 2950
$ws.Range("M80").Value = -1952

$ws.Range("H83").Value = 3360
$ws.Range("I83").Value = 2950
$ws.Range("K83").Value = 14750
$ws.Range("M83").Value = -9758

$ws.Range("H132").Value = 11368875
$ws.Range("I132").Value = 11368875
$ws.Range("K132").Value = 34106625
$ws.Range("M132").Value = -34104095

$ws.Range("H136").Value = 22084.15
$ws.Range("J136").Value = 22084.15
$ws.Range("L136").Value = 66252.45000000001
$ws.Range("N136").Value = -71352.45000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5000
$ws.Range("I7").Value = 5000
$ws.Range("K7").Value = 5000
$ws.Range("M7").Value = -4888

$ws.Range("H61").Value = 2938.0557
$ws.Range("J61").Value = 3248.5
$ws.Range("L61").Value = 3248.5
$ws.Range("N61").Value = -3652.5

$ws.Range("H82").Value = 984.8
$ws.Range("I82").Value = 991.6667
$ws.Range("J82").Value = 974.5
$ws.Range("K82").Value = 991.6667
$ws.Range("L82").Value = 974.5
$ws.Range("M82").Value = -630.6667
$ws.Range("N82").Value = -1696.5

$ws.Range("H85").Value = 984.8
$ws.Range("I85").Value = 991.6667
$ws.Range("J85").Value = 974.5
$ws.Range("K85").Value = 991.6667
$ws.Range("L85").Value = 974.5
$ws.Range("M85").Value = 256.3333
$ws.Range("N85").Value = -3470.5

$ws.Range("H106").Value = 16773
$ws.Range("J106").Value = 16773
$ws.Range("L106").Value = 16773
$ws.Range("N106").Value = -19297

$ws.Range("H113").Value = 2938.0557
$ws.Range("J113").Value = 3248.5
$ws.Range("L113").Value = 3248.5
$ws.Range("N113").Value = -7588.5

$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1040
$ws.Range("I100").Value = 1151.5714
$ws.Range("J100").Value = 844.75
$ws.Range("K100").Value = 2303.1428
$ws.Range("L100").Value = 1689.5
$ws.Range("M100").Value = -1762.1428
$ws.Range("N100").Value = -2771.5

$ws.Range("H105").Value = 29999.5
$ws.Range("I105").Value = 19999
$ws.Range("J105").Value = 40000
$ws.Range("K105").Value = 19999
$ws.Range("L105").Value = 40000
$ws.Range("M105").Value = -16505
$ws.Range("N105").Value = -46988

$ws.Range("H132").Value = 21748134
$ws.Range("I132").Value = 25005060
$ws.Range("J132").Value = 35306
$ws.Range("K132").Value = 75015180
$ws.Range("L132").Value = 105918
$ws.Range("M132").Value = -75012650
$ws.Range("N132").Value = -110978
